$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Role" column header and values to use the new PI/Sub I terminology.
$ws.Range("A1").Value = "Role (PI/Sub I)"
$ws.Range("A2").Value = "PI"
$ws.Range("A3").Value = "Sub I"
$ws.Range("A4").Value = "Sub I"
$ws.Range("A5").Value = "Sub I"

# Move the active selection to A5, matching the saved view state.
$ws.Range("A5").Select()
